$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'58.038.45"
$ws.Cells.Item(2, 5).Value = "'  -0.39%  "
$ws.Cells.Item(3, 4).Value = "'2.454.30"
$ws.Cells.Item(3, 5).Value = "'  -2.75%  "
$ws.Cells.Item(4, 5).Value = "'  -0.03%  "
$ws.Cells.Item(5, 4).Value = "'524.96"
$ws.Cells.Item(5, 5).Value = "'  +0.55%  "
$ws.Cells.Item(6, 4).Value = "'131.69"
$ws.Cells.Item(6, 5).Value = "'  -1.01%  "
$ws.Cells.Item(7, 5).Value = "'  +0.01%  "
$ws.Cells.Item(8, 4).Value = "'0.565"
$ws.Cells.Item(8, 5).Value = "'  +0.42%  "
$ws.Cells.Item(9, 4).Value = "'2.461.01"
$ws.Cells.Item(9, 5).Value = "'  -2.45%  "
$ws.Cells.Item(10, 5).Value = "'  +0.41%  "
$ws.Cells.Item(11, 4).Value = "'0.151"
$ws.Cells.Item(11, 5).Value = "'  -1.75%  "
$ws.Cells.Item(12, 4).Value = "'4.98"
$ws.Cells.Item(12, 5).Value = "'  -3.49%  "
$ws.Cells.Item(13, 5).Value = "'  -1.96%  "
$ws.Cells.Item(14, 4).Value = "'2.890.56"
$ws.Cells.Item(14, 5).Value = "'  -2.56%  "
$ws.Cells.Item(15, 4).Value = "'57.988.20"
$ws.Cells.Item(15, 5).Value = "'  -0.53%  "
$ws.Cells.Item(16, 4).Value = "'21.82"
$ws.Cells.Item(16, 5).Value = "'  -1.41%  "
$ws.Cells.Item(17, 5).Value = "'  -1.23%  "
$ws.Cells.Item(18, 4).Value = "'2.461.80"
$ws.Cells.Item(18, 5).Value = "'  -2.30%  "
$ws.Cells.Item(19, 4).Value = "'10.34"
$ws.Cells.Item(19, 5).Value = "'  -3.01%  "
$ws.Cells.Item(20, 4).Value = "'4.13"
$ws.Cells.Item(20, 5).Value = "'  -0.71%  "
$ws.Cells.Item(21, 4).Value = "'311.71"
$ws.Cells.Item(21, 5).Value = "'  -3.14%  "
$ws.Cells.Item(22, 5).Value = "'  -1.00%  "
$ws.Cells.Item(23, 5).Value = "'  +0.02%  "
$ws.Cells.Item(24, 4).Value = "'65.06"
$ws.Cells.Item(24, 5).Value = "'  +0.91%  "
$ws.Cells.Item(25, 5).Value = "'  -0.65%  "
$ws.Cells.Item(26, 4).Value = "'2.587.10"
$ws.Cells.Item(26, 5).Value = "'  -1.44%  "
$ws.Cells.Item(27, 4).Value = "'0.998"
$ws.Cells.Item(27, 5).Value = "'  -0.14%  "
$ws.Cells.Item(28, 5).Value = "'  -1.38%  "
$ws.Cells.Item(29, 5).Value = "'  -2.09%  "
$ws.Cells.Item(30, 4).Value = "'173.39"
$ws.Cells.Item(30, 5).Value = "'  +2.68%  "
$ws.Cells.Item(31, 5).Value = "'  -1.70%  "
$ws.Cells.Item(32, 5).Value = "'  -1.14%  "
$ws.Cells.Item(33, 5).Value = "'  -0.87%  "
$ws.Cells.Item(34, 5).Value = "'  -3.89%  "
$ws.Cells.Item(35, 4).Value = "'0.998"
$ws.Cells.Item(35, 5).Value = "'  +0.01%  "
$ws.Cells.Item(36, 4).Value = "'0.998"
$ws.Cells.Item(36, 5).Value = "'  +0.05%  "
$ws.Cells.Item(37, 4).Value = "'17.83"
$ws.Cells.Item(37, 5).Value = "'  -1.84%  "
$ws.Cells.Item(38, 5).Value = "'  -5.04%  "
$ws.Cells.Item(39, 5).Value = "'  -3.07%  "
$ws.Cells.Item(40, 4).Value = "'0.819"
$ws.Cells.Item(40, 5).Value = "'  +6.56%  "
$ws.Cells.Item(41, 4).Value = "'36.25"
$ws.Cells.Item(41, 5).Value = "'  -0.54%  "
$ws.Cells.Item(42, 5).Value = "'  -1.94%  "
$ws.Cells.Item(43, 5).Value = "'  -0.95%  "
$ws.Cells.Item(44, 4).Value = "'262.72"
$ws.Cells.Item(44, 5).Value = "'  -5.01%  "
$ws.Cells.Item(45, 4).Value = "'0.588"
$ws.Cells.Item(45, 5).Value = "'  -1.73%  "
$ws.Cells.Item(46, 4).Value = "'4.82"
$ws.Cells.Item(46, 5).Value = "'  -3.21%  "
$ws.Cells.Item(47, 5).Value = "'  +0.59%  "
$ws.Cells.Item(48, 4).Value = "'122.29"
$ws.Cells.Item(48, 5).Value = "'  -5.94%  "
$ws.Cells.Item(49, 5).Value = "'  -0.89%  "
$ws.Cells.Item(50, 5).Value = "'  -0.81%  "
$ws.Cells.Item(51, 4).Value = "'17.03"
$ws.Cells.Item(51, 5).Value = "'  -3.75%  "
